$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header row values (shared strings): date + fund codes + stock tickers
$headers = @(
  "date",
  "RBF263.CF(RBC U.S. Equity Fund Series A)",
  "RBF590.CF(RBC U.S. Dividend Fund Series A)",
  "RBF557.CF(RBC U.S. Index Fund Series Dz)",
  "INA48603.CF(IA US Daq Index (Ia) SRP7575Myed+)",
  "INA36081.CF(iA Global Innovators Fid SRP75100)",
  "AAPL(Apple Inc. (AAPL))",
  "META(Meta Platforms, Inc. (META))",
  "NVDA(NVIDIA Corporation (NVDA))",
  "AMZN(Amazon.com, Inc. (AMZN))",
  "MSFT(Microsoft Corporation (MSFT))",
  "SHOP(Shopify Inc. (SHOP))",
  "TSLA(Tesla, Inc. (TSLA))",
  "GOOG(Alphabet Inc. (GOOG))",
  "AVGO(Broadcom Inc. (AVGO))"
)

for ($i = 0; $i -lt $headers.Length; $i++) {
  $ws.Cells.Item(1, $i + 1).Value = $headers[$i]
}

# Re-style the whole header row: Calibri 11 (theme text colour), centered,
# thin black box border around every cell (replaces the old bold 宋体 style).
$full = $ws.Range("A1:O1")
$full.ClearFormats()
$full.HorizontalAlignment = -4108
$full.Borders.LineStyle = 1
$full.Borders.Color = 0
$full.Font.Name = "Calibri"

# The data-header cells (everything but the date column) also get a
# #,##0.00 number format applied.
$dataRng = $ws.Range("B1:O1")
$dataRng.NumberFormat = "#,##0.00"

# Row is a bit taller now and the new selection spans the header row.
$full.RowHeight = 15
[void]$ws.Range("A1:P1").Select()
